# Apply the row permutation for rows 36-39 on the Artfynd sheet.
# Each row's full record (A:AY) moves to a new row number:
#   old row 36 -> new row 37
#   old row 37 -> new row 39
#   old row 38 -> new row 36
#   old row 39 -> new row 38
#
# Because this is a 4-cycle, a direct sequential copy would overwrite
# source data before it is read. Stage each source row in a scratch
# area first (using real Copy/PasteSpecial so cell types/formatting are
# preserved exactly, not re-interpreted the way a plain .Value write
# would be), then move the staged rows into their final homes, then
# wipe the scratch area. Every destination is cleared before pasting so
# that cells blank in the source do not leave stale data behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch36 = "A1000:AY1000"
$scratch37 = "A1001:AY1001"
$scratch38 = "A1002:AY1002"
$scratch39 = "A1003:AY1003"

# Stage all four source rows first.
$ws.Range($scratch36).ClearContents()
$ws.Range("A36:AY36").Copy()
$ws.Range("A1000").PasteSpecial()

$ws.Range($scratch37).ClearContents()
$ws.Range("A37:AY37").Copy()
$ws.Range("A1001").PasteSpecial()

$ws.Range($scratch38).ClearContents()
$ws.Range("A38:AY38").Copy()
$ws.Range("A1002").PasteSpecial()

$ws.Range($scratch39).ClearContents()
$ws.Range("A39:AY39").Copy()
$ws.Range("A1003").PasteSpecial()

# Now write the staged rows into their destination rows.
$ws.Range("A37:AY37").ClearContents()
$ws.Range($scratch36).Copy()
$ws.Range("A37").PasteSpecial()

$ws.Range("A39:AY39").ClearContents()
$ws.Range($scratch37).Copy()
$ws.Range("A39").PasteSpecial()

$ws.Range("A36:AY36").ClearContents()
$ws.Range($scratch38).Copy()
$ws.Range("A36").PasteSpecial()

$ws.Range("A38:AY38").ClearContents()
$ws.Range($scratch39).Copy()
$ws.Range("A38").PasteSpecial()

# Clean up the scratch area.
$ws.Range("A1000:AY1003").ClearContents()
